# Updated cryptos list on Mon Nov 27 16:00:19 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (price / volume%) to be treated as plain text so
# Excel does not silently reinterpret values like "54.49" or "0.0780" as
# numbers (which would change their stored representation).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.878.08"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").Value = "2.008.08"
$ws.Range("E3").Value = "  -2.81%  "

$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "224.78"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("E6").Value = "  -3.64%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "54.49"
$ws.Range("E8").Value = "  -4.76%  "

$ws.Range("D9").Value = "0.377"
$ws.Range("E9").Value = "  -2.56%  "

$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("E11").Value = "  -4.18%  "

$ws.Range("D12").Value = "2.308.58"
$ws.Range("E12").Value = "  -2.67%  "

$ws.Range("D13").Value = "14.18"
$ws.Range("E13").Value = "  -4.31%  "

$ws.Range("D14").Value = "20.23"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "0.736"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").Value = "5.10"
$ws.Range("E16").Value = "  -3.40%  "

$ws.Range("D17").Value = "2.022.05"
$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("D18").Value = "36.763.74"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("E19").Value = "  +4.86%  "

$ws.Range("D20").Value = "68.61"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").Value = "0.0₃0816"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").Value = "225.02"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -7.44%  "

$ws.Range("D26").Value = "165.49"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("E27").Value = "  -4.01%  "

$ws.Range("E28").Value = "  -4.79%  "

$ws.Range("D29").Value = "18.56"
$ws.Range("E29").Value = "  -3.64%  "

$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").Value = "4.39"
$ws.Range("E34").Value = "  -3.82%  "

$ws.Range("E35").Value = "  -4.96%  "

$ws.Range("E36").Value = "  +1.68%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  -5.04%  "

$ws.Range("E39").Value = "  +0.25%  "

$ws.Range("D40").Value = "1.477.39"
$ws.Range("E40").Value = "  +0.76%  "

$ws.Range("E41").Value = "  -5.31%  "

$ws.Range("D42").Value = "16.75"
$ws.Range("E42").Value = "  +1.84%  "

$ws.Range("D43").Value = "94.78"
$ws.Range("E43").Value = "  -4.05%  "

$ws.Range("D44").Value = "0.0920"
$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("E45").Value = "  -5.10%  "

$ws.Range("E46").Value = "  -5.93%  "

$ws.Range("D47").Value = "7.23"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("D50").Value = "2.198.20"
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("E51").Value = "  -12.03%  "
